# Juno: check in to OLPRODLOC.
#
# For each of the nine "section intro" paragraphs, split the leading
# sentence run into two runs: a short bold "label" (e.g. "目的", "协作",
# "领导：") followed by the (non-bold) remainder of the sentence, whose
# wording is also refreshed to the new Chinese text.

$d = $word.ActiveDocument

function Split-LabelRun($ParaIndex, $OldText, $BoldText, $RestText) {
    $p = $d.Paragraphs($ParaIndex)
    $pStart = $p.Range.Start

    # Sanity-check: the paragraph's leading run must still hold the text
    # we expect to replace, so we never clobber the wrong run.
    $check = $d.Range($pStart, $pStart + $OldText.Length)
    if ($check.Text -ne $OldText) {
        throw "Paragraph $ParaIndex did not contain expected text. Found: $($check.Text)"
    }

    # Rewrite the whole run's text to label + rest, then re-derive the
    # range (Word collapses/expands Start/End to the freshly set text).
    $check.Text = $BoldText + $RestText

    # Bold just the label portion; Word automatically splits the run
    # into two at this boundary, leaving the remainder with the run's
    # original (non-bold) character formatting.
    $boldRange = $d.Range($pStart, $pStart + $BoldText.Length)
    $boldRange.Font.Bold = 1
}

Split-LabelRun 3 `
    "本文件概述了平面设计学院所有设计团队成员的核心职责。" `
    "目的" `
    "：本文档概述了图形设计研究所所有设计团队成员的核心职责。"

Split-LabelRun 5 `
    "与其他设计师、开发人员和利益干系人合作，创造符合项目要求的高质量设计。" `
    "协作" `
    "：与其他设计人员、开发人员和利益干系人协作，创建满足项目要求的高质量设计。"

Split-LabelRun 13 `
    "创造具有视觉吸引力、便于用户使用、易于访问且响应速度快的设计。" `
    "设计" `
    "：创建具有视觉吸引力的设计，这些设计对用户友好、可访问和响应性强。"

Split-LabelRun 22 `
    "与团队成员、利益干系人和客户进行有效沟通，确保满足项目要求。" `
    "沟通" `
    "：与团队成员、利益干系人和客户有效沟通，以确保满足项目要求。"

Split-LabelRun 30 `
    "开展研究，确定用户需求、偏好和行为，为设计决策提供依据。" `
    "研究" `
    "：进行研究以确定用户需求、偏好和行为，以告知设计决策。"

Split-LabelRun 39 `
    "进行可用性测试，确保设计满足用户需求，便于所有用户使用。" `
    "测试" `
    "：进行可用性测试，以确保设计满足用户需求，可供所有用户访问。"

Split-LabelRun 48 `
    "创建并维护设计文档，包括设计规范、风格指南和设计模式。" `
    "文档" `
    "：创建和维护设计文档，包括设计规范、样式指南和设计模式。"

Split-LabelRun 56 `
    "了解最新的设计趋势、工具和技术，提高设计质量和效率。" `
    "专业开发" `
    "：随时了解最新的设计趋势、工具和技术，以提高设计质量和效率。"

Split-LabelRun 61 `
    "领导设计团队，为初级设计师提供指导。" `
    "领导：" `
    " 领导设计团队，为初级设计师提供指导。"

Write-Output "done"
